$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.203.69'
$ws.Range("E2").Value = '  +2.86%  '
$ws.Range("D3").Value = '2.652.55'
$ws.Range("E3").Value = '  +2.57%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.01'
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.43'
$ws.Range("E6").Value = '  +3.75%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +1.15%  '
$ws.Range("E9").Value = '  +7.56%  '
$ws.Range("E10").Value = '  +4.08%  '
$ws.Range("E11").Value = '  +1.59%  '
$ws.Range("E12").Value = '  +1.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.05'
$ws.Range("E13").Value = '  +5.28%  '
$ws.Range("E14").Value = '  +18.81%  '
$ws.Range("D15").Value = '3.129.29'
$ws.Range("E15").Value = '  +2.61%  '
$ws.Range("D16").Value = '65.049.57'
$ws.Range("E16").Value = '  +2.94%  '
$ws.Range("D17").Value = '2.689.05'
$ws.Range("E17").Value = '  +4.92%  '
$ws.Range("E18").Value = '  +3.17%  '
$ws.Range("E19").Value = '  +1.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '355.10'
$ws.Range("E20").Value = '  +2.78%  '
$ws.Range("E21").Value = '  +6.08%  '
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.71'
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("E25").Value = '  +2.93%  '
$ws.Range("E26").Value = '  -1.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.18'
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.165'
$ws.Range("E28").Value = '  +1.84%  '
$ws.Range("D29").Value = '0.0₃0948'
$ws.Range("E29").Value = '  +11.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '522.32'
$ws.Range("E31").Value = '  -6.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.11'
$ws.Range("E32").Value = '  +3.65%  '
$ws.Range("E33").Value = '  +1.80%  '
$ws.Range("E34").Value = '  +7.16%  '
$ws.Range("E35").Value = '  +4.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.428'
$ws.Range("E36").Value = '  +3.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '165.09'
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '20.31'
$ws.Range("E38").Value = '  +3.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.02'
$ws.Range("E39").Value = '  +5.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.15'
$ws.Range("E42").Value = '  +6.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '165.61'
$ws.Range("E43").Value = '  -0.51%  '
$ws.Range("E44").Value = '  +2.75%  '
$ws.Range("E45").Value = '  +5.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.07'
$ws.Range("E46").Value = '  +0.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.21'
$ws.Range("E47").Value = '  +4.53%  '
$ws.Range("E48").Value = '  +3.55%  '
$ws.Range("E49").Value = '  +1.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0988'
$ws.Range("E50").Value = '  +2.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.53'
$ws.Range("E51").Value = '  +1.49%  '
